$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows scraped from the diff: dates 2024-09-02 .. 2024-09-27 appended as rows 674-693.
# Each entry: Date, Colgate(C), ITC(D), Marico(E), Emami(F), VarunBev(G), Basket Value(H), daily return(I), NAV(J)
$rows = @(
    @("2024-09-02", 650.9500122070312, 1470.050048828125, 608.5800170898438, 1316.800048828125, 811.2000122070312, 17987.74047851562, 0, 179.8587255790135),
    @("2024-09-03", 640.0499877929688, 1460.75, 599.9400024414062, 1341.949951171875, 822.3499755859375, 17977.66967773438, -0.0005598702512568748, 179.7580280291328),
    @("2024-09-04", 645.5999755859375, 1475.300048828125, 609, 1327.75, 824.2999877929688, 18052.54992675781, 0.004165181047695957, 180.506752760651),
    @("2024-09-05", 643.9000244140625, 1457.699951171875, 602.1799926757812, 1254.800048828125, 835.4000244140625, 17792.94024658203, -0.01438077619112318, 177.9109255482136),
    @("2024-09-06", 665.25, 1443.449951171875, 597.2999877929688, 1256.849975585938, 832.7000122070312, 17880.34979248047, 0.004912597057432855, 178.7849302375469),
    @("2024-09-09", 676, 1492.050048828125, 610.3400268554688, 1225.25, 827.5999755859375, 18025.32012939453, 0.008107802061849448, 180.2344830635545),
    @("2024-09-10", 680, 1503.050048828125, 608, 1246, 824.75, 18130.15014648438, 0.005815709032478913, 181.2826743746714),
    @("2024-09-11", 680.4500122070312, 1499.949951171875, 627.6599731445312, 1229, 814, 18088.97985839844, -0.002270818926114677, 180.8710142467247),
    @("2024-09-12", 686.0999755859375, 1513.449951171875, 645.5999755859375, 1224.849975585938, 809.7000122070312, 18193.19958496094, 0.00576150382046627, 181.9131032863188),
    @("2024-09-13", 681.9500122070312, 1491.300048828125, 646.6500244140625, 1229.300048828125, 788.0499877929688, 18027.60040283203, -0.009102257211853799, 180.2572834300002),
    @("2024-09-16", 695.2000122070312, 1456.349975585938, 621.0499877929688, 1219.699951171875, 751.9500122070312, 17765.49987792969, -0.01453884704817227, 177.6365503568924),
    @("2024-09-17", 692, 1459.400024414062, 649.6500244140625, 1222.949951171875, 746.75, 17827, 0.003461772676980224, 178.2514877133509),
    @("2024-09-18", 695.2999877929688, 1432.150024414062, 646.7000122070312, 1224.550048828125, 744.5999755859375, 17755.70007324219, -0.003999547133999692, 177.5385624865358),
    @("2024-09-19", 697, 1444.849975585938, 652.1500244140625, 1197.849975585938, 747.2000122070312, 17752.34997558594, -0.0001886773060161447, 177.5050649888519),
    @("2024-09-20", 709, 1456.599975585938, 654.4500122070312, 1206.300048828125, 747.5499877929688, 17905.25006103516, 0.008612949027001824, 179.0339070656355),
    @("2024-09-23", 702.5, 1449.300048828125, 654.0999755859375, 1190, 763.75, 17852.70007324219, -0.00293489270542646, 178.5084617577646),
    @("2024-09-24", 705.0999755859375, 1446.349975585938, 646.8499755859375, 1194.699951171875, 760.9500122070312, 17843.19958496094, -0.0005321597429113499, 178.4134667406481),
    @("2024-09-25", 689.2000122070312, 1429.550048828125, 633.2999877929688, 1175.349975585938, 742.5499877929688, 17509.20007324219, -0.0187185885652627, 175.0738184622277),
    @("2024-09-26", 693.5999755859375, 1422.300048828125, 626.8499755859375, 1165, 742.25, 17466.64990234375, -0.002430160756656341, 174.6483609390828),
    @("2024-09-27", 692.4500122070312, 1388.650024414062, 608.5499877929688, 1165.550048828125, 751.6500244140625, 17342.00036621094, -0.007136430673868747, 173.4019950189362)
)

$startRow = 674
$r = $startRow
foreach ($row in $rows) {
    # Column A: literal text date (leading apostrophe forces text, matching the inlineStr dates already in the sheet)
    $ws.Cells.Item($r, 1).Formula = "'" + $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]   # C: Colgate-Palmolive
    $ws.Cells.Item($r, 4).Value = $row[2]   # D: ITC
    $ws.Cells.Item($r, 5).Value = $row[3]   # E: Marico
    $ws.Cells.Item($r, 6).Value = $row[4]   # F: Emami
    $ws.Cells.Item($r, 7).Value = $row[5]   # G: Varun Beverages
    $ws.Cells.Item($r, 8).Value = $row[6]   # H: Basket Value
    # Column I inherits the sheet-wide "0.0%" column style; reset to Normal so it stays
    # unformatted like the rest of the recent rows (matching rows 594-673 already in the sheet).
    $ws.Cells.Item($r, 9).Style = "Normal"
    $ws.Cells.Item($r, 9).Value = $row[7]   # I: daily return
    $ws.Cells.Item($r, 10).Value = $row[8]  # J: NAV
    $r = $r + 1
}

Write-Host "Inserted $($rows.Count) rows starting at row $startRow"